$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.878.19'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.14%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.490.26'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.03%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '149.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.99%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.541'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.33%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.490.14'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.98%  '
$ws.Range("E10").Value = '  +1.12%  '
$ws.Range("E11").Value = '  +1.04%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.361'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.52%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.29'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.67%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.36'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.89%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000182'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.63%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.945.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.95%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '63.673.26'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.98%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.499.50'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.38'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '331.43'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.11'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +20.37%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.30'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.68%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '628.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000105'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.74%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.69'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.64%  '
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.00%  '
$ws.Range("B30").Value = 'WrappedeETH'
$ws.Range("C30").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.623.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.67%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.45'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.24%  '
$ws.Range("B32").Value = 'Binance-PegBSC-USD'
$ws.Range("C32").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.14%  '
$ws.Range("E33").Value = '  -2.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.92'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.15%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.999'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.386'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.52'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.45%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.84'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.33%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '148.61'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.74'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +13.73%  '
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '150.44'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.05%  '
$ws.Range("E46").Value = '  +2.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.44'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0548'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.609'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.31%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0237'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0923'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.52%  '